$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates for sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 568
$ws1.Range("F5").Value = 282
$ws1.Range("F6").Value = 1086
$ws1.Range("F7").Value = 1430
$ws1.Range("F9").Value = 106
$ws1.Range("F10").Value = 747
$ws1.Range("F15").Value = 1343
$ws1.Range("F16").Value = 108
$ws1.Range("F17").Value = 97
$ws1.Range("F20").Value = 649
$ws1.Range("F23").Value = 214
$ws1.Range("F24").Value = 21
$ws1.Range("F25").Value = 5803
$ws1.Range("F26").Value = 61
$ws1.Range("F28").Value = 94
$ws1.Range("F30").Value = 14407
$ws1.Range("F31").Value = 1434
$ws1.Range("F32").Value = 204
$ws1.Range("F33").Value = 101
$ws1.Range("F34").Value = 86
$ws1.Range("F35").Value = 7009
$ws1.Range("F36").Value = 608
$ws1.Range("F37").Value = 4195
$ws1.Range("F39").Value = 358
$ws1.Range("F40").Value = 111

# Column F ("想去人数") updates for sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 568
$ws4.Range("F5").Value = 282
$ws4.Range("F6").Value = 1086
$ws4.Range("F7").Value = 1430
$ws4.Range("F9").Value = 106
$ws4.Range("F10").Value = 747
$ws4.Range("F15").Value = 1343
$ws4.Range("F16").Value = 108
$ws4.Range("F17").Value = 97
$ws4.Range("F21").Value = 649
$ws4.Range("F25").Value = 214
$ws4.Range("F26").Value = 21
$ws4.Range("F28").Value = 5804
$ws4.Range("F29").Value = 61
$ws4.Range("F31").Value = 94
$ws4.Range("F33").Value = 14407
$ws4.Range("F34").Value = 1434
$ws4.Range("F35").Value = 204
$ws4.Range("F36").Value = 101
$ws4.Range("F37").Value = 86
$ws4.Range("F38").Value = 7020
$ws4.Range("F39").Value = 608
$ws4.Range("F40").Value = 4195
$ws4.Range("F42").Value = 358
$ws4.Range("F43").Value = 111
